$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A87").Value = 1342412855899868928
$ws.Range("B87").Value = "Eh eh ada yg introvert gak nih disini? `r`nNih mimin kasih jurusan cocok buat kmuu ☺️`r`nCek 👇`r`nhttps://t.co/6bnYERsLTR`r`n#jurusankuliah #snmptn2021 #infosbmptn #utbk2021 #beasiswa #infomenarik #kuliahonline #FYP #masukkampus #kampusfavorit #ugm #ui #undip #ltmpt #anaksma #pejuang21 https://t.co/ViMg8DUJh5"
$ws.Range("C87").Value = "kampungsoal"
$ws.Range("D87").Value = "Fri Dec 25 10:12:20 +0000 2020"

$ws.Range("A88").Value = 1342076779947692032
$ws.Range("B88").Value = "Ketika lihat omongan dosen di grup kelas `r`nDosen : `" anak-anak mata kuliah dimulai 5 menit lagi`"`r`nRequest buat alasan yng paling manjur? `r`n#kuliahonline`r`n#kuliah"
$ws.Range("C88").Value = "cimolse"
$ws.Range("D88").Value = "Thu Dec 24 11:56:53 +0000 2020"

$ws.Range("A89").Value = 1341969450556222976
$ws.Range("B89").Value = "Ternyata serba online itu sulit ya.`r`nDari kuliah online, sekolah online, dan chattingan 7 hari sekilas centang biru padahal online. Syukur syukur adaptasi nya gak lewat online juga ... Leee `r`n#online #kuliahonline #serbaonline"
$ws.Range("C89").Value = "DedekGu58931839"
$ws.Range("D89").Value = "Thu Dec 24 04:50:24 +0000 2020"

$ws.Range("A90").Value = 1341796168225837056
$ws.Range("B90").Value = "uang segitu, ngga layak sih sama apa yang didapet,denger temen yang udah kerja juga malah jadi pengen ikutan kerja.inget tujuan awal heh udah ngorbanin apa aja buat bisa masuk kuliah :)`r`n#kuliahonline #coronaasu #salahesopojal #cowoktetepsalah"
$ws.Range("C90").Value = "coklatpahitmen"
$ws.Range("D90").Value = "Wed Dec 23 17:21:50 +0000 2020"

$ws.Range("A91").Value = 1341705698497809920
$ws.Range("B91").Value = "Kalimat paling horror pas kulon bukan`r`n`"ayo siapa yg mau bertanya`"`r`nTapi ...`r`n`"Ayo buka kameranya saya hitung sampe 3`"`r`nlangsung kelayapan nyari kerudung 😬😬`r`n#kuliahonline"
$ws.Range("C91").Value = "Piyaaaakkk"
$ws.Range("D91").Value = "Wed Dec 23 11:22:20 +0000 2020"

$ws.Range("A92").Value = 1341549080900562944
$ws.Range("B92").Value = "⏩PAKET DESIGN #ADOBE ILLUSTRATOR🤩 - Pembuatan Design Menggunakan Adobe Illustrator 18 Module https://t.co/V7y8EtrLeU`r`n#PaketDesign #adobeillustrator #pakaimasker #jagajarak #cucitangan #dirumahaja #kerjadarirumah #belajardirumah #KuliahOnline #IndonesiaSehat #ai #jagakesehatan"
$ws.Range("C92").Value = "PotekantropusX"
$ws.Range("D92").Value = "Wed Dec 23 01:00:00 +0000 2020"

$ws.Range("A93").Value = 1341386047612023040
$ws.Range("B93").Value = "Yg penting selalu berubah dgn menjadi baik dan bermanfaat. Ilmu saat kuliah untuk modal awal cara berpikir dan bertindak`r`n#ITBlogAwards #pendidikan #kuliah #KuliahOnline #Training #Menkes #menag #sandiagauno #ReshuffleKabinet #ilmufardukifayah #YaqutCholilQoumas"
$ws.Range("C93").Value = "MKompetensi"
$ws.Range("D93").Value = "Tue Dec 22 14:12:10 +0000 2020"

$ws.Range("A94").Value = 1341320794790452992
$ws.Range("B94").Value = "Jangan Ragu Kuliah dengan konsep yang fleksibel dan biaya terjangkau , UMT solusinya .`r`n#umtindonesia #umt #maba2021 #mabacovid #KuliahOnline https://t.co/QKzdFiLa25"
$ws.Range("C94").Value = "UMT_Indonesia"
$ws.Range("D94").Value = "Tue Dec 22 09:52:52 +0000 2020"

$ws.Range("A95").Value = 1341316291802480896
$ws.Range("B95").Value = "Pengertian Kebijakan Mutu https://t.co/sUErkXK7nz #mutu #kebijakanmutu #Quality #qualitycontrol #kualitas #manajemenmutu #KuliahOnline #teknikindustri #industrial #America #qualitypolicy #seventools #lean #Germany #England #Jakarta #University #manajemen #surabaya #sarjana #uii https://t.co/IatX9ngVJO"
$ws.Range("C95").Value = "charif_noor"
$ws.Range("D95").Value = "Tue Dec 22 09:34:59 +0000 2020"

$ws.Range("A96").Value = 1341186694091919104
$ws.Range("B96").Value = "⏩PAKET DESIGN #ADOBE ILLUSTRATOR🤩 - Pembuatan Design Menggunakan Adobe Illustrator 18 Module https://t.co/V7y8EtrLeU`r`n#PaketDesign #adobeillustrator #pakaimasker #jagajarak #cucitangan #dirumahaja #kerjadarirumah #belajardirumah #KuliahOnline #IndonesiaSehat #ai #jagakesehatan"
$ws.Range("C96").Value = "PotekantropusX"
$ws.Range("D96").Value = "Tue Dec 22 01:00:00 +0000 2020"

$ws.Range("A97").Value = 1341175988797853952
$ws.Range("B97").Value = "Bangun langsung cuci muka ❎`r`nBangun langsung presentasi ✔`r`n#kuliah #kuliahonline #ceritapagi #daring #lelucon #mendesak #presentasi #sayamautobat"
$ws.Range("C97").Value = "orhandy_"
$ws.Range("D97").Value = "Tue Dec 22 00:17:28 +0000 2020"

# Reset row heights back to default (undoes the automatic "wrap text" row-
# height bump that setting a multi-line cell value triggers) so the new rows
# match the plain, unstyled rows already present in the sheet.
$ws.Rows("87:97").AutoFit()

# Update the selected/active cell to match the post-edit view state
$ws.Range("I88").Select()
